$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("E2") "720"
Set-TextValue $ws.Range("F2") "8"
Set-TextValue $ws.Range("G2") "8"
Set-TextValue $ws.Range("J3") "3"
Set-TextValue $ws.Range("E7") "810"
Set-TextValue $ws.Range("F7") "9"
Set-TextValue $ws.Range("G7") "9"
Set-TextValue $ws.Range("L7") "4"
Set-TextValue $ws.Range("E8") "563"
Set-TextValue $ws.Range("F8") "9"
Set-TextValue $ws.Range("G8") "5"
Set-TextValue $ws.Range("E9") "187"
Set-TextValue $ws.Range("F9") "4"
Set-TextValue $ws.Range("G9") "2"
Set-TextValue $ws.Range("I9") "1"
Set-TextValue $ws.Range("E10") "237"
Set-TextValue $ws.Range("F10") "5"
Set-TextValue $ws.Range("G10") "3"
Set-TextValue $ws.Range("E11") "16"
Set-TextValue $ws.Range("F11") "1"
Set-TextValue $ws.Range("H11") "1"
Set-TextValue $ws.Range("J11") "3"
Set-TextValue $ws.Range("E12") "810"
Set-TextValue $ws.Range("F12") "9"
Set-TextValue $ws.Range("G12") "9"
Set-TextValue $ws.Range("J14") "5"
Set-TextValue $ws.Range("F15") "9"
Set-TextValue $ws.Range("H15") "3"
Set-TextValue $ws.Range("J15") "3"
Set-TextValue $ws.Range("E16") "305"
Set-TextValue $ws.Range("F16") "9"
Set-TextValue $ws.Range("H16") "7"
Set-TextValue $ws.Range("J16") "7"
Set-TextValue $ws.Range("E17") "522"
Set-TextValue $ws.Range("F17") "8"
Set-TextValue $ws.Range("G17") "7"
Set-TextValue $ws.Range("I17") "7"
Set-TextValue $ws.Range("E18") "512"
Set-TextValue $ws.Range("F18") "8"
Set-TextValue $ws.Range("G18") "6"
Set-TextValue $ws.Range("I18") "4"
Set-TextValue $ws.Range("L18") "5"
Set-TextValue $ws.Range("J19") "7"
Set-TextValue $ws.Range("J21") "2"
Set-TextValue $ws.Range("F22") "6"
Set-TextValue $ws.Range("H22") "4"
Set-TextValue $ws.Range("J22") "7"
Set-TextValue $ws.Range("E24") "760"
Set-TextValue $ws.Range("F24") "9"
Set-TextValue $ws.Range("G24") "9"
Set-TextValue $ws.Range("I24") "4"
Set-TextValue $ws.Range("E25") "722"
Set-TextValue $ws.Range("F25") "9"
Set-TextValue $ws.Range("G25") "8"
Set-TextValue $ws.Range("E26") "571"
Set-TextValue $ws.Range("F26") "9"
Set-TextValue $ws.Range("G26") "7"
Set-TextValue $ws.Range("I26") "5"
Set-TextValue $ws.Range("E27") "236"
Set-TextValue $ws.Range("F27") "7"
Set-TextValue $ws.Range("H27") "5"
Set-TextValue $ws.Range("J27") "6"
